$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("mmWave(BR)")
$ws.Range("A143:A162").NumberFormat = "@"
$ws.Range("B143:C162").NumberFormat = "@"
$ws.Cells.Item(143, 1).Value = "2026-02-01"
$ws.Cells.Item(143, 2).Value = "21:09:03"
$ws.Cells.Item(143, 3).Value = "21:00"
$ws.Cells.Item(143, 4).Value = "Bedroom"
$ws.Cells.Item(143, 5).Value = 1
$ws.Cells.Item(143, 6).Value = "Occupied"
$ws.Cells.Item(144, 1).Value = "2026-02-01"
$ws.Cells.Item(144, 2).Value = "21:09:04"
$ws.Cells.Item(144, 3).Value = "21:00"
$ws.Cells.Item(144, 4).Value = "Bedroom"
$ws.Cells.Item(144, 5).Value = 2
$ws.Cells.Item(144, 6).Value = "Occupied"
$ws.Cells.Item(145, 1).Value = "2026-02-01"
$ws.Cells.Item(145, 2).Value = "21:09:06"
$ws.Cells.Item(145, 3).Value = "21:00"
$ws.Cells.Item(145, 4).Value = "Bedroom"
$ws.Cells.Item(145, 5).Value = 1
$ws.Cells.Item(145, 6).Value = "Occupied"
$ws.Cells.Item(146, 1).Value = "2026-02-01"
$ws.Cells.Item(146, 2).Value = "21:09:07"
$ws.Cells.Item(146, 3).Value = "21:00"
$ws.Cells.Item(146, 4).Value = "Bedroom"
$ws.Cells.Item(146, 5).Value = 4
$ws.Cells.Item(146, 6).Value = "Occupied"
$ws.Cells.Item(147, 1).Value = "2026-02-01"
$ws.Cells.Item(147, 2).Value = "21:09:08"
$ws.Cells.Item(147, 3).Value = "21:00"
$ws.Cells.Item(147, 4).Value = "Bedroom"
$ws.Cells.Item(147, 5).Value = 2
$ws.Cells.Item(147, 6).Value = "Occupied"
$ws.Cells.Item(148, 1).Value = "2026-02-01"
$ws.Cells.Item(148, 2).Value = "21:09:10"
$ws.Cells.Item(148, 3).Value = "21:00"
$ws.Cells.Item(148, 4).Value = "Bedroom"
$ws.Cells.Item(148, 5).Value = 1
$ws.Cells.Item(148, 6).Value = "Occupied"
$ws.Cells.Item(149, 1).Value = "2026-02-01"
$ws.Cells.Item(149, 2).Value = "21:09:11"
$ws.Cells.Item(149, 3).Value = "21:00"
$ws.Cells.Item(149, 4).Value = "Bedroom"
$ws.Cells.Item(149, 5).Value = 3
$ws.Cells.Item(149, 6).Value = "Occupied"
$ws.Cells.Item(150, 1).Value = "2026-02-01"
$ws.Cells.Item(150, 2).Value = "21:09:12"
$ws.Cells.Item(150, 3).Value = "21:00"
$ws.Cells.Item(150, 4).Value = "Bedroom"
$ws.Cells.Item(150, 5).Value = 2
$ws.Cells.Item(150, 6).Value = "Occupied"
$ws.Cells.Item(151, 1).Value = "2026-02-01"
$ws.Cells.Item(151, 2).Value = "21:09:14"
$ws.Cells.Item(151, 3).Value = "21:00"
$ws.Cells.Item(151, 4).Value = "Bedroom"
$ws.Cells.Item(151, 5).Value = 1
$ws.Cells.Item(151, 6).Value = "Occupied"
$ws.Cells.Item(152, 1).Value = "2026-02-01"
$ws.Cells.Item(152, 2).Value = "21:09:43"
$ws.Cells.Item(152, 3).Value = "21:00"
$ws.Cells.Item(152, 4).Value = "Bedroom"
$ws.Cells.Item(152, 5).Value = 23
$ws.Cells.Item(152, 6).Value = "Occupied"
$ws.Cells.Item(153, 1).Value = "2026-02-01"
$ws.Cells.Item(153, 2).Value = "21:09:44"
$ws.Cells.Item(153, 3).Value = "21:00"
$ws.Cells.Item(153, 4).Value = "Bedroom"
$ws.Cells.Item(153, 5).Value = 2
$ws.Cells.Item(153, 6).Value = "Occupied"
$ws.Cells.Item(154, 1).Value = "2026-02-01"
$ws.Cells.Item(154, 2).Value = "21:09:46"
$ws.Cells.Item(154, 3).Value = "21:00"
$ws.Cells.Item(154, 4).Value = "Bedroom"
$ws.Cells.Item(154, 5).Value = 3
$ws.Cells.Item(154, 6).Value = "Occupied"
$ws.Cells.Item(155, 1).Value = "2026-02-01"
$ws.Cells.Item(155, 2).Value = "21:09:47"
$ws.Cells.Item(155, 3).Value = "21:00"
$ws.Cells.Item(155, 4).Value = "Bedroom"
$ws.Cells.Item(155, 5).Value = 2
$ws.Cells.Item(155, 6).Value = "Occupied"
$ws.Cells.Item(156, 1).Value = "2026-02-01"
$ws.Cells.Item(156, 2).Value = "21:09:49"
$ws.Cells.Item(156, 3).Value = "21:00"
$ws.Cells.Item(156, 4).Value = "Bedroom"
$ws.Cells.Item(156, 5).Value = 3
$ws.Cells.Item(156, 6).Value = "Occupied"
$ws.Cells.Item(157, 1).Value = "2026-02-01"
$ws.Cells.Item(157, 2).Value = "21:09:50"
$ws.Cells.Item(157, 3).Value = "21:00"
$ws.Cells.Item(157, 4).Value = "Bedroom"
$ws.Cells.Item(157, 5).Value = 2
$ws.Cells.Item(157, 6).Value = "Occupied"
$ws.Cells.Item(158, 1).Value = "2026-02-01"
$ws.Cells.Item(158, 2).Value = "21:09:51"
$ws.Cells.Item(158, 3).Value = "21:00"
$ws.Cells.Item(158, 4).Value = "Bedroom"
$ws.Cells.Item(158, 5).Value = 3
$ws.Cells.Item(158, 6).Value = "Occupied"
$ws.Cells.Item(159, 1).Value = "2026-02-01"
$ws.Cells.Item(159, 2).Value = "21:09:53"
$ws.Cells.Item(159, 3).Value = "21:00"
$ws.Cells.Item(159, 4).Value = "Bedroom"
$ws.Cells.Item(159, 5).Value = 2
$ws.Cells.Item(159, 6).Value = "Occupied"
$ws.Cells.Item(160, 1).Value = "2026-02-01"
$ws.Cells.Item(160, 2).Value = "21:09:55"
$ws.Cells.Item(160, 3).Value = "21:00"
$ws.Cells.Item(160, 4).Value = "Bedroom"
$ws.Cells.Item(160, 5).Value = 4
$ws.Cells.Item(160, 6).Value = "Occupied"
$ws.Cells.Item(161, 1).Value = "2026-02-01"
$ws.Cells.Item(161, 2).Value = "21:09:56"
$ws.Cells.Item(161, 3).Value = "21:00"
$ws.Cells.Item(161, 4).Value = "Bedroom"
$ws.Cells.Item(161, 5).Value = 2
$ws.Cells.Item(161, 6).Value = "Occupied"
$ws.Cells.Item(162, 1).Value = "2026-02-01"
$ws.Cells.Item(162, 2).Value = "21:10:01"
$ws.Cells.Item(162, 3).Value = "21:00"
$ws.Cells.Item(162, 4).Value = "Bedroom"
$ws.Cells.Item(162, 5).Value = 3
$ws.Cells.Item(162, 6).Value = "Occupied"
$ws = $wb.Worksheets.Item("mmWave(HR)")
$ws.Range("A144:A163").NumberFormat = "@"
$ws.Range("B144:C163").NumberFormat = "@"
$ws.Cells.Item(144, 1).Value = "2026-02-01"
$ws.Cells.Item(144, 2).Value = "21:09:03"
$ws.Cells.Item(144, 3).Value = "21:00"
$ws.Cells.Item(144, 4).Value = "Bedroom"
$ws.Cells.Item(144, 5).Value = 49
$ws.Cells.Item(144, 6).Value = "Occupied"
$ws.Cells.Item(145, 1).Value = "2026-02-01"
$ws.Cells.Item(145, 2).Value = "21:09:04"
$ws.Cells.Item(145, 3).Value = "21:00"
$ws.Cells.Item(145, 4).Value = "Bedroom"
$ws.Cells.Item(145, 5).Value = 50
$ws.Cells.Item(145, 6).Value = "Occupied"
$ws.Cells.Item(146, 1).Value = "2026-02-01"
$ws.Cells.Item(146, 2).Value = "21:09:05"
$ws.Cells.Item(146, 3).Value = "21:00"
$ws.Cells.Item(146, 4).Value = "Bedroom"
$ws.Cells.Item(146, 5).Value = 49
$ws.Cells.Item(146, 6).Value = "Occupied"
$ws.Cells.Item(147, 1).Value = "2026-02-01"
$ws.Cells.Item(147, 2).Value = "21:09:07"
$ws.Cells.Item(147, 3).Value = "21:00"
$ws.Cells.Item(147, 4).Value = "Bedroom"
$ws.Cells.Item(147, 5).Value = 52
$ws.Cells.Item(147, 6).Value = "Occupied"
$ws.Cells.Item(148, 1).Value = "2026-02-01"
$ws.Cells.Item(148, 2).Value = "21:09:08"
$ws.Cells.Item(148, 3).Value = "21:00"
$ws.Cells.Item(148, 4).Value = "Bedroom"
$ws.Cells.Item(148, 5).Value = 50
$ws.Cells.Item(148, 6).Value = "Occupied"
$ws.Cells.Item(149, 1).Value = "2026-02-01"
$ws.Cells.Item(149, 2).Value = "21:09:09"
$ws.Cells.Item(149, 3).Value = "21:00"
$ws.Cells.Item(149, 4).Value = "Bedroom"
$ws.Cells.Item(149, 5).Value = 49
$ws.Cells.Item(149, 6).Value = "Occupied"
$ws.Cells.Item(150, 1).Value = "2026-02-01"
$ws.Cells.Item(150, 2).Value = "21:09:11"
$ws.Cells.Item(150, 3).Value = "21:00"
$ws.Cells.Item(150, 4).Value = "Bedroom"
$ws.Cells.Item(150, 5).Value = 51
$ws.Cells.Item(150, 6).Value = "Occupied"
$ws.Cells.Item(151, 1).Value = "2026-02-01"
$ws.Cells.Item(151, 2).Value = "21:09:12"
$ws.Cells.Item(151, 3).Value = "21:00"
$ws.Cells.Item(151, 4).Value = "Bedroom"
$ws.Cells.Item(151, 5).Value = 50
$ws.Cells.Item(151, 6).Value = "Occupied"
$ws.Cells.Item(152, 1).Value = "2026-02-01"
$ws.Cells.Item(152, 2).Value = "21:09:13"
$ws.Cells.Item(152, 3).Value = "21:00"
$ws.Cells.Item(152, 4).Value = "Bedroom"
$ws.Cells.Item(152, 5).Value = 49
$ws.Cells.Item(152, 6).Value = "Occupied"
$ws.Cells.Item(153, 1).Value = "2026-02-01"
$ws.Cells.Item(153, 2).Value = "21:09:42"
$ws.Cells.Item(153, 3).Value = "21:00"
$ws.Cells.Item(153, 4).Value = "Bedroom"
$ws.Cells.Item(153, 5).Value = 71
$ws.Cells.Item(153, 6).Value = "Occupied"
$ws.Cells.Item(154, 1).Value = "2026-02-01"
$ws.Cells.Item(154, 2).Value = "21:09:44"
$ws.Cells.Item(154, 3).Value = "21:00"
$ws.Cells.Item(154, 4).Value = "Bedroom"
$ws.Cells.Item(154, 5).Value = 50
$ws.Cells.Item(154, 6).Value = "Occupied"
$ws.Cells.Item(155, 1).Value = "2026-02-01"
$ws.Cells.Item(155, 2).Value = "21:09:45"
$ws.Cells.Item(155, 3).Value = "21:00"
$ws.Cells.Item(155, 4).Value = "Bedroom"
$ws.Cells.Item(155, 5).Value = 51
$ws.Cells.Item(155, 6).Value = "Occupied"
$ws.Cells.Item(156, 1).Value = "2026-02-01"
$ws.Cells.Item(156, 2).Value = "21:09:47"
$ws.Cells.Item(156, 3).Value = "21:00"
$ws.Cells.Item(156, 4).Value = "Bedroom"
$ws.Cells.Item(156, 5).Value = 50
$ws.Cells.Item(156, 6).Value = "Occupied"
$ws.Cells.Item(157, 1).Value = "2026-02-01"
$ws.Cells.Item(157, 2).Value = "21:09:48"
$ws.Cells.Item(157, 3).Value = "21:00"
$ws.Cells.Item(157, 4).Value = "Bedroom"
$ws.Cells.Item(157, 5).Value = 51
$ws.Cells.Item(157, 6).Value = "Occupied"
$ws.Cells.Item(158, 1).Value = "2026-02-01"
$ws.Cells.Item(158, 2).Value = "21:09:49"
$ws.Cells.Item(158, 3).Value = "21:00"
$ws.Cells.Item(158, 4).Value = "Bedroom"
$ws.Cells.Item(158, 5).Value = 50
$ws.Cells.Item(158, 6).Value = "Occupied"
$ws.Cells.Item(159, 1).Value = "2026-02-01"
$ws.Cells.Item(159, 2).Value = "21:09:51"
$ws.Cells.Item(159, 3).Value = "21:00"
$ws.Cells.Item(159, 4).Value = "Bedroom"
$ws.Cells.Item(159, 5).Value = 51
$ws.Cells.Item(159, 6).Value = "Occupied"
$ws.Cells.Item(160, 1).Value = "2026-02-01"
$ws.Cells.Item(160, 2).Value = "21:09:52"
$ws.Cells.Item(160, 3).Value = "21:00"
$ws.Cells.Item(160, 4).Value = "Bedroom"
$ws.Cells.Item(160, 5).Value = 50
$ws.Cells.Item(160, 6).Value = "Occupied"
$ws.Cells.Item(161, 1).Value = "2026-02-01"
$ws.Cells.Item(161, 2).Value = "21:09:54"
$ws.Cells.Item(161, 3).Value = "21:00"
$ws.Cells.Item(161, 4).Value = "Bedroom"
$ws.Cells.Item(161, 5).Value = 52
$ws.Cells.Item(161, 6).Value = "Occupied"
$ws.Cells.Item(162, 1).Value = "2026-02-01"
$ws.Cells.Item(162, 2).Value = "21:09:56"
$ws.Cells.Item(162, 3).Value = "21:00"
$ws.Cells.Item(162, 4).Value = "Bedroom"
$ws.Cells.Item(162, 5).Value = 50
$ws.Cells.Item(162, 6).Value = "Occupied"
$ws.Cells.Item(163, 1).Value = "2026-02-01"
$ws.Cells.Item(163, 2).Value = "21:10:00"
$ws.Cells.Item(163, 3).Value = "21:00"
$ws.Cells.Item(163, 4).Value = "Bedroom"
$ws.Cells.Item(163, 5).Value = 51
$ws.Cells.Item(163, 6).Value = "Occupied"
$ws = $wb.Worksheets.Item("mmWave(InBed)")
$ws.Range("A147:A166").NumberFormat = "@"
$ws.Range("B147:C166").NumberFormat = "@"
$ws.Cells.Item(147, 1).Value = "2026-02-01"
$ws.Cells.Item(147, 2).Value = "21:09:02"
$ws.Cells.Item(147, 3).Value = "21:00"
$ws.Cells.Item(147, 4).Value = "Bedroom"
$ws.Cells.Item(147, 5).Value = "In Bed"
$ws.Cells.Item(147, 6).Value = "Occupied"
$ws.Cells.Item(148, 1).Value = "2026-02-01"
$ws.Cells.Item(148, 2).Value = "21:09:03"
$ws.Cells.Item(148, 3).Value = "21:00"
$ws.Cells.Item(148, 4).Value = "Bedroom"
$ws.Cells.Item(148, 5).Value = "In Bed"
$ws.Cells.Item(148, 6).Value = "Occupied"
$ws.Cells.Item(149, 1).Value = "2026-02-01"
$ws.Cells.Item(149, 2).Value = "21:09:05"
$ws.Cells.Item(149, 3).Value = "21:00"
$ws.Cells.Item(149, 4).Value = "Bedroom"
$ws.Cells.Item(149, 5).Value = "In Bed"
$ws.Cells.Item(149, 6).Value = "Occupied"
$ws.Cells.Item(150, 1).Value = "2026-02-01"
$ws.Cells.Item(150, 2).Value = "21:09:06"
$ws.Cells.Item(150, 3).Value = "21:00"
$ws.Cells.Item(150, 4).Value = "Bedroom"
$ws.Cells.Item(150, 5).Value = "In Bed"
$ws.Cells.Item(150, 6).Value = "Occupied"
$ws.Cells.Item(151, 1).Value = "2026-02-01"
$ws.Cells.Item(151, 2).Value = "21:09:08"
$ws.Cells.Item(151, 3).Value = "21:00"
$ws.Cells.Item(151, 4).Value = "Bedroom"
$ws.Cells.Item(151, 5).Value = "In Bed"
$ws.Cells.Item(151, 6).Value = "Occupied"
$ws.Cells.Item(152, 1).Value = "2026-02-01"
$ws.Cells.Item(152, 2).Value = "21:09:09"
$ws.Cells.Item(152, 3).Value = "21:00"
$ws.Cells.Item(152, 4).Value = "Bedroom"
$ws.Cells.Item(152, 5).Value = "In Bed"
$ws.Cells.Item(152, 6).Value = "Occupied"
$ws.Cells.Item(153, 1).Value = "2026-02-01"
$ws.Cells.Item(153, 2).Value = "21:09:10"
$ws.Cells.Item(153, 3).Value = "21:00"
$ws.Cells.Item(153, 4).Value = "Bedroom"
$ws.Cells.Item(153, 5).Value = "In Bed"
$ws.Cells.Item(153, 6).Value = "Occupied"
$ws.Cells.Item(154, 1).Value = "2026-02-01"
$ws.Cells.Item(154, 2).Value = "21:09:12"
$ws.Cells.Item(154, 3).Value = "21:00"
$ws.Cells.Item(154, 4).Value = "Bedroom"
$ws.Cells.Item(154, 5).Value = "In Bed"
$ws.Cells.Item(154, 6).Value = "Occupied"
$ws.Cells.Item(155, 1).Value = "2026-02-01"
$ws.Cells.Item(155, 2).Value = "21:09:13"
$ws.Cells.Item(155, 3).Value = "21:00"
$ws.Cells.Item(155, 4).Value = "Bedroom"
$ws.Cells.Item(155, 5).Value = "In Bed"
$ws.Cells.Item(155, 6).Value = "Occupied"
$ws.Cells.Item(156, 1).Value = "2026-02-01"
$ws.Cells.Item(156, 2).Value = "21:09:42"
$ws.Cells.Item(156, 3).Value = "21:00"
$ws.Cells.Item(156, 4).Value = "Bedroom"
$ws.Cells.Item(156, 5).Value = "In Bed"
$ws.Cells.Item(156, 6).Value = "Occupied"
$ws.Cells.Item(157, 1).Value = "2026-02-01"
$ws.Cells.Item(157, 2).Value = "21:09:43"
$ws.Cells.Item(157, 3).Value = "21:00"
$ws.Cells.Item(157, 4).Value = "Bedroom"
$ws.Cells.Item(157, 5).Value = "In Bed"
$ws.Cells.Item(157, 6).Value = "Occupied"
$ws.Cells.Item(158, 1).Value = "2026-02-01"
$ws.Cells.Item(158, 2).Value = "21:09:45"
$ws.Cells.Item(158, 3).Value = "21:00"
$ws.Cells.Item(158, 4).Value = "Bedroom"
$ws.Cells.Item(158, 5).Value = "In Bed"
$ws.Cells.Item(158, 6).Value = "Occupied"
$ws.Cells.Item(159, 1).Value = "2026-02-01"
$ws.Cells.Item(159, 2).Value = "21:09:46"
$ws.Cells.Item(159, 3).Value = "21:00"
$ws.Cells.Item(159, 4).Value = "Bedroom"
$ws.Cells.Item(159, 5).Value = "In Bed"
$ws.Cells.Item(159, 6).Value = "Occupied"
$ws.Cells.Item(160, 1).Value = "2026-02-01"
$ws.Cells.Item(160, 2).Value = "21:09:48"
$ws.Cells.Item(160, 3).Value = "21:00"
$ws.Cells.Item(160, 4).Value = "Bedroom"
$ws.Cells.Item(160, 5).Value = "In Bed"
$ws.Cells.Item(160, 6).Value = "Occupied"
$ws.Cells.Item(161, 1).Value = "2026-02-01"
$ws.Cells.Item(161, 2).Value = "21:09:49"
$ws.Cells.Item(161, 3).Value = "21:00"
$ws.Cells.Item(161, 4).Value = "Bedroom"
$ws.Cells.Item(161, 5).Value = "In Bed"
$ws.Cells.Item(161, 6).Value = "Occupied"
$ws.Cells.Item(162, 1).Value = "2026-02-01"
$ws.Cells.Item(162, 2).Value = "21:09:50"
$ws.Cells.Item(162, 3).Value = "21:00"
$ws.Cells.Item(162, 4).Value = "Bedroom"
$ws.Cells.Item(162, 5).Value = "In Bed"
$ws.Cells.Item(162, 6).Value = "Occupied"
$ws.Cells.Item(163, 1).Value = "2026-02-01"
$ws.Cells.Item(163, 2).Value = "21:09:52"
$ws.Cells.Item(163, 3).Value = "21:00"
$ws.Cells.Item(163, 4).Value = "Bedroom"
$ws.Cells.Item(163, 5).Value = "In Bed"
$ws.Cells.Item(163, 6).Value = "Occupied"
$ws.Cells.Item(164, 1).Value = "2026-02-01"
$ws.Cells.Item(164, 2).Value = "21:09:54"
$ws.Cells.Item(164, 3).Value = "21:00"
$ws.Cells.Item(164, 4).Value = "Bedroom"
$ws.Cells.Item(164, 5).Value = "In Bed"
$ws.Cells.Item(164, 6).Value = "Occupied"
$ws.Cells.Item(165, 1).Value = "2026-02-01"
$ws.Cells.Item(165, 2).Value = "21:09:55"
$ws.Cells.Item(165, 3).Value = "21:00"
$ws.Cells.Item(165, 4).Value = "Bedroom"
$ws.Cells.Item(165, 5).Value = "In Bed"
$ws.Cells.Item(165, 6).Value = "Occupied"
$ws.Cells.Item(166, 1).Value = "2026-02-01"
$ws.Cells.Item(166, 2).Value = "21:10:00"
$ws.Cells.Item(166, 3).Value = "21:00"
$ws.Cells.Item(166, 4).Value = "Bedroom"
$ws.Cells.Item(166, 5).Value = "In Bed"
$ws.Cells.Item(166, 6).Value = "Occupied"
